$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Sheet 1')
$ws.Range('E5').Value = 'Based on official disease reports to the WOAH'
$ws.Range('E6').Value = 'AHS is a disease listed in the World Organisation for Animal Health ({ref009:WOAH}) Terrestrial Animal Health Code and must be reported to the WOAH. The map to the right displays outbreak points reported to the WOAH early warning system since 2005.'
$ws.Range('E7').Value = 'As described in the WOAH {ref005:Terrestrial Animal Health Code}, the WOAH early warning system includes immediate notifications and follow-up reports on:'
$ws.Range('E14').Value = 'Countries are coloured according to the available information regarding their stable disease situation (disease status legend). This information is provided by countries through the WOAH monitoring system, which is a different reporting channel.<br>Immediate notifications (points) and disease status (country/region colours) are reported to the WOAH in different spatial and temporal scales, and therefore are displayed in the map as layers which can be filtered independently.'
$ws.Range('E17').Value = 'For more up to date reports, visit the original data source: {ref001:WOAH-WAHIS}.'
$ws.Range('E31').Value = 'A summary of the disease in animal hosts is given in the {ref008:WOAH Technical disease card}.'
$ws.Range('E44').Value = 'Humans are not susceptible to AHSV and therefore there is no direct impact on public health ({ref008:WOAH Technical disease card}).'
$ws.Range('E53').Value = 'Refer to the {ref008:WOAH Technical disease card} for a key summary of the virus characteristics. '
$ws.Range('E65').Value = 'Refer to the {ref008:WOAH Technical disease card} for a key summary of the disease transmission and epidemiological parameters.'
$ws.Range('E77').Value = 'WOAH prescribed tests for international trade include enzyme-linked immunosorbent assay (ELISA) and complement fixation (CF) tests ({ref010:WOAH, Terrestrial Manual})'
$ws.Range('E97').Value = 'Bites of a biological vector of the susceptible animal hosts (horses, mules, donkeys and zebras) are required for transmission of AHSV. The most significant vector* seems to be <i>Culicoides imicola</i>, but other species, such as <i>C. variipennis</i>, which is common in many parts of the United States and C. bolitinos, present in Africa, should also be considered as potential vectors ({ref034:Boinas et al., 2009}). Occasional transmission can occur (this was only experimental) from mosquitoes (<i>Culex</i>, <i>Anopheles</i> and <i>Aedes spp.</i>); ticks (<i>Hyalomma dromedari</i>, <i>Rhipicephalus sanguineus</i>) and, possibly, biting flies (<i>Stomoxys</i> and <i>Tabanus</i>) ({ref008:WOAH technical disease card}).'
$ws.Range('E113').Value = 'Several attenuated (monovalent and polyvalent) live vaccines for use in horses, mules and donkeys, are currently commercially available. Inactivated or recombinant vaccines are not commercialised, but some killed vaccines are licensed for use in some areas (e.g. Egypt) ({ref008: WOAH, Technical Disease Card}). Currently, no vaccines have been authorised for use in the European Union by the European Medicine Agency ({ref035:EMA}).'
$ws.Range('E118').Value = 'There are no specific curative treatments for AHSV infections ({ref008: WOAH, Technical Disease Card}).'
$ws.Range('E130').Value = 'Geographical distribution data has been kindly provided by the World Organisation of Animal Health (WOAH). {ref001:WOAH-WAHIS} (WOAH World Animal Health Information System) is the original source of these data.'
$ws = $wb.Worksheets.Item('References')
$ws.Range('C2').Value = 'WOAH-WAHIS (WOAH World Animal Health Information System)'
$ws.Range('C5').Value = 'WOAH (World Organisation for Animal Health). Terrestrial Animal Health Code 2021. WOAH, Paris, France'
$ws.Range('C8').Value = 'WOAH (World Organisation for Animal Health) Technical Disease Card: African Horse Sickness. 2021.'
$ws.Range('C9').Value = 'WOAH (World Organisation for Animal Health), 2021. African Swine fever. Chapter 12.1. WOAH Terrestrial Animal Health Code 2021. WOAH, Paris, France'
$ws.Range('C10').Value = 'WOAH (World Organisation for Animal Health), 2019. African Swine fever. Chapter 3.5.1. WOAH Terrestrial Manual 2019. WOAH, Paris, France'
